$d = $word.ActiveDocument

# The paragraph currently reads (across 3 differently-formatted runs):
#   "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (" +
#   "http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/" +
#   ")."
# It needs to become a single run of plain text with the year bumped to 2022:
#   "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$oldText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Clear the whole (multi-run) span, then type the replacement text fresh so it
    # lands in a single run with no leftover run-level formatting overrides.
    $rng.Text = ""
    $rng.InsertAfter($newText)
}
